$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared string rich-text runs) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- C20: numeric 2 -> text "0" (requires quote-prefix-free style 14) ---
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null

# --- Data table value updates (rows 15-30) ---
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -83.333333333333
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -23.529411764705
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5.084745762711
$ws.Range("M16").Value = -15.068493150684
$ws.Range("N16").Value = -78.245614035087
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 6.666666666666
$ws.Range("I17").Value = 94
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = 13.253012048192
$ws.Range("L17").Value = 28.767123287671
$ws.Range("M17").Value = 38.235294117647
$ws.Range("N17").Value = -58.407079646017
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -58.823529411764
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = -37.681159420289
$ws.Range("L18").Value = -15.686274509803
$ws.Range("M18").Value = -4.444444444444
$ws.Range("N18").Value = -86.604361370716
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -10.714285714285
$ws.Range("I19").Value = 129
$ws.Range("J19").Value = 124
$ws.Range("K19").Value = 4.032258064516
$ws.Range("L19").Value = 16.216216216216
$ws.Range("M19").Value = 27.722772277227
$ws.Range("N19").Value = -14.569536423841
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 3
$ws.Range("E20").NumberFormat = "#,##0.0;"-"#,##0.0"
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 25
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = 5.263157894736
$ws.Range("L20").Value = 11.111111111111
$ws.Range("M20").Value = 185.714285714286
$ws.Range("N20").Value = -65.517241379310
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 19.047619047619
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -18.518518518518
$ws.Range("I21").Value = 353
$ws.Range("J21").Value = 360
$ws.Range("K21").Value = -1.944444444444
$ws.Range("L21").Value = 11.356466876971
$ws.Range("M21").Value = 18.060200668896
$ws.Range("N21").Value = -66.854460093896
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;"-"#,##0.0"
$ws.Range("E22").Value = 0
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -60
$ws.Range("L22").Value = -14.285714285714
$ws.Range("M22").Value = -25
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = -11.428571428571
$ws.Range("L23").Value = 6.896551724137
$ws.Range("M23").Value = 181.818181818182
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 38.095238095238
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 25
$ws.Range("I24").Value = 449
$ws.Range("J24").Value = 399
$ws.Range("K24").Value = 12.531328320802
$ws.Range("L24").Value = 25.418994413407
$ws.Range("M24").Value = 31.671554252199
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 225
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 41.176470588235
$ws.Range("I25").Value = 160
$ws.Range("J25").Value = 161
$ws.Range("K25").Value = -0.621118012422
$ws.Range("L25").Value = 5.263157894736
$ws.Range("M25").Value = -9.090909090909
$ws.Range("L26").Value = -50
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -34.782608695652
$ws.Range("N28").Value = -83.333333333333
$ws.Range("N29").Value = -83.870967741935
$ws.Range("L30").Value = 50
